# Update cryptocurrency price/volume data per the scraped source refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.361.90"
$ws.Range("E2").Value = "  -2.58%  "
$ws.Range("D3").Value = "2.540.81"
$ws.Range("E3").Value = "  -4.21%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "514.62"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.68%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.34"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.38%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.558"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.14%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.48"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -7.60%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0993"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.20%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.325"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.70%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.130"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.01%  "
$ws.Range("D13").Value = "2.988.40"
$ws.Range("D14").Value = "57.360.82"
$ws.Range("E14").Value = "  -2.61%  "
$ws.Range("E15").Value = "  -4.95%  "
$ws.Range("E16").Value = "  -2.95%  "
$ws.Range("D17").Value = "2.584.67"
$ws.Range("E17").Value = "  -2.12%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "332.68"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.86%  "
$ws.Range("E19").Value = "  -1.83%  "
$ws.Range("E20").Value = "  -2.38%  "
$ws.Range("E21").Value = "  -3.68%  "
$ws.Range("E22").Value = "  -0.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "64.61"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.12%  "
$ws.Range("E24").Value = "  -0.14%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.400"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.36%  "
$ws.Range("D27").Value = "2.667.30"
$ws.Range("E27").Value = "  -3.79%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.94"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.99%  "
$ws.Range("E29").Value = "  -6.16%  "
$ws.Range("E30").Value = "  +0.04%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.24"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -6.73%  "
$ws.Range("E32").Value = "  -2.14%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "148.84"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.53%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "18.46"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.99%  "
$ws.Range("E35").Value = "  -4.38%  "
$ws.Range("E36").Value = "  -5.32%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.837"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.78%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "35.73"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.79%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.823"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.33%  "
$ws.Range("E40").Value = "  -4.64%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.09%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.47"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.13%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.62"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.44%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0951"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.62%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.573"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -7.07%  "
$ws.Range("B46").Value = "Hedera"
$ws.Range("C46").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0519"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.58%  "
$ws.Range("B47").Value = "Bittensor"
$ws.Range("C47").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "257.53"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.50%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "18.37"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -7.76%  "
$ws.Range("D49").Value = "1.964.88"
$ws.Range("E49").Value = "  -3.45%  "
$ws.Range("E50").Value = "  -2.86%  "
$ws.Range("E51").Value = "  -4.68%  "
